# se modif datos de cuenta
$wb = $excel.ActiveWorkbook

$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsHogar  = $wb.Worksheets.Item("DatosHogar")
$wsMotor  = $wb.Worksheets.Item("DatosMotor")
$wsAP     = $wb.Worksheets.Item("DatosAP")

# --- DatosCuenta (sheet1) ---
$wsCuenta.Range("A2").Value = "SmokeDos"
$wsCuenta.Range("B2").Value = "SmokeLastDos"
$wsCuenta.Range("C2").Value = 20111102
$wsCuenta.Range("D2").Value = 102

# --- DatosHogar (sheet2) ---
$wsHogar.Range("A2").Value = 622

# --- DatosMotor (sheet3) ---
$wsMotor.Range("A2").Value = "SMA003"
$wsMotor.Range("B2").Value = "ABC12SSMA003"
$wsMotor.Range("C2").Value = "ZAZ123SSMA003"

# --- DatosAP (sheet4) ---
$wsAP.Range("A2").Value = 21200102

# --- Selections / active cells / active sheet ---
$wsHogar.Range("A3").Select()
$wsAP.Range("H2").Select()
$wsMotor.Range("A2:C2").Select()

$wsCuenta.Activate()
$wsCuenta.Range("B2").Select()
